# "open file functionality and calculation part added"
#
# Adds a numeric helper column (D) used for the open/calculation logic,
# updates a few measured values, and moves the selection to D6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New helper cells in column D (calculation inputs, default 0) ---
$ws.Range("D5").Value  = 0
$ws.Range("D9").Value  = 0
$ws.Range("D10").Value = 0
$ws.Range("D11").Value = 0
$ws.Range("D12").Value = 0
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0

# --- Row 12: "Linear refractive index" value becomes a real number (1.8)
#     instead of the text placeholder "1", keeping a text-like display format.
$ws.Range("C12").Value = 1.8
$ws.Range("C12").NumberFormat = "@"

# --- Row 13: "Transmitance" measured value changed from 20 % to 85 % ---
$ws.Range("C13").Value = "85 %"

# --- Row 14: "Thickness" measured value changed from 1 mm to 2 mm ---
$ws.Range("C14").Value = "2 mm"

# --- Selection moved to D6 (as left by the editing session) ---
$ws.Range("D6").Select()
